$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.600.93"
$ws.Range("E2").Value = "  -5.69%  "
$ws.Range("D3").Value = "2.210.30"
$ws.Range("E3").Value = "  -7.31%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "486.83"
$ws.Range("E5").Value = "  -4.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "125.32"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("E8").Value = "  -4.61%  "
$ws.Range("D9").Value = "2.237.98"
$ws.Range("E9").Value = "  -6.68%  "
$ws.Range("E10").Value = "  -6.86%  "
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.318"
$ws.Range("E12").Value = "  -3.17%  "
$ws.Range("E13").Value = "  -4.38%  "
$ws.Range("D14").Value = "2.605.95"
$ws.Range("E14").Value = "  -7.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.26"
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("D16").Value = "53.531.81"
$ws.Range("E16").Value = "  -5.58%  "
$ws.Range("E17").Value = "  -4.28%  "
$ws.Range("D18").Value = "2.226.97"
$ws.Range("E18").Value = "  -6.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.67"
$ws.Range("E19").Value = "  -4.70%  "
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "295.75"
$ws.Range("E21").Value = "  -4.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.17"
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.54"
$ws.Range("E24").Value = "  -4.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.996"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.367"
$ws.Range("E26").Value = "  -1.49%  "
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").Value = "2.311.95"
$ws.Range("E28").Value = "  -7.41%  "
$ws.Range("E29").Value = "  -3.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "163.00"
$ws.Range("E30").Value = "  -6.01%  "
$ws.Range("E31").Value = "  -4.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.997"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "0.0₃0671"
$ws.Range("E33").Value = "  -6.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.79"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.992"
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("E36").Value = "  -2.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.33"
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.841"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("E40").Value = "  -4.93%  "
$ws.Range("E41").Value = "  -3.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.368"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("E43").Value = "  -1.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.30"
$ws.Range("E44").Value = "  -3.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "126.77"
$ws.Range("E45").Value = "  -2.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.83"
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("E48").Value = "  -6.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "233.99"
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0472"
$ws.Range("E50").Value = "  -2.66%  "
$ws.Range("E51").Value = "  -3.65%  "
